$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 37333.332
$ws.Range("I9").Value = 100001
$ws.Range("J9").Value = 5999.5
$ws.Range("K9").Value = 100001
$ws.Range("L9").Value = 5999.5
$ws.Range("M9").Value = -99832
$ws.Range("N9").Value = -6337.5
$ws.Range("H41").Value = 90.71429000000001
$ws.Range("I41").Value = 70.545456
$ws.Range("J41").Value = 164.66667
$ws.Range("K41").Value = 70.545456
$ws.Range("L41").Value = 164.66667
$ws.Range("M41").Value = 369.454544
$ws.Range("N41").Value = -1044.66667
$ws.Range("H92").Value = 546.7619
$ws.Range("I92").Value = 568.1
$ws.Range("J92").Value = 120
$ws.Range("K92").Value = 568.1
$ws.Range("L92").Value = 120
$ws.Range("M92").Value = 679.9
$ws.Range("N92").Value = -2616
$ws.Range("H106").Value = 1978.2222
$ws.Range("I106").Value = 1850.5
$ws.Range("K106").Value = 1850.5
$ws.Range("M106").Value = -1219.5
$ws.Range("H135").Value = 41667760
$ws.Range("I135").Value = 1262.25
$ws.Range("K135").Value = 11360.25
$ws.Range("M135").Value = -8825.25
$ws.Range("H137").Value = 5693.846
$ws.Range("I137").Value = 2845.7666
$ws.Range("K137").Value = 8537.299800000001
$ws.Range("M137").Value = -5987.299800000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5195.347
$ws.Range("I32").Value = 4608.952
$ws.Range("J32").Value = 8713.714
$ws.Range("K32").Value = 4608.952
$ws.Range("L32").Value = 8713.714
$ws.Range("M32").Value = -4321.952
$ws.Range("N32").Value = -9287.714
$ws.Range("H45").Value = 3919.1667
$ws.Range("I45").Value = 4394.6
$ws.Range("J45").Value = 3324.875
$ws.Range("K45").Value = 4394.6
$ws.Range("L45").Value = 3324.875
$ws.Range("M45").Value = -4017.6
$ws.Range("N45").Value = -4078.875
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2262.2334
$ws.Range("I20").Value = 1614.2106
$ws.Range("J20").Value = 3381.5454
$ws.Range("K20").Value = 1614.2106
$ws.Range("L20").Value = 3381.5454
$ws.Range("M20").Value = -1367.2106
$ws.Range("N20").Value = -3875.5454
$ws.Range("H80").Value = 950
$ws.Range("J80").Value = 920.8
$ws.Range("L80").Value = 920.8
$ws.Range("N80").Value = -2916.8
$ws.Range("H83").Value = 950
$ws.Range("J83").Value = 920.8
$ws.Range("L83").Value = 4604
$ws.Range("N83").Value = -14588
$ws.Range("H94").Value = 1532.3914
$ws.Range("I94").Value = 1314.4546
$ws.Range("J94").Value = 1732.1666
$ws.Range("K94").Value = 1314.4546
$ws.Range("L94").Value = 1732.1666
$ws.Range("M94").Value = -863.4546
$ws.Range("N94").Value = -2634.1666
$ws.Range("H102").Value = 52220.57
$ws.Range("J102").Value = 93333.336
$ws.Range("L102").Value = 93333.336
$ws.Range("N102").Value = -99823.336
$ws.Range("H134").Value = 5283.16
$ws.Range("I134").Value = 5312.174
$ws.Range("J134").Value = 4949.5
$ws.Range("K134").Value = 15936.522
$ws.Range("L134").Value = 14848.5
$ws.Range("M134").Value = -13401.522
$ws.Range("N134").Value = -19918.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15631084
$ws.Range("I31").Value = 3505
$ws.Range("K31").Value = 3505
$ws.Range("M31").Value = -3210
$ws.Range("H34").Value = 15631084
$ws.Range("I34").Value = 3505
$ws.Range("K34").Value = 3505
$ws.Range("M34").Value = -3303
$ws.Range("H58").Value = 1987.7858
$ws.Range("I58").Value = 1146.4445
$ws.Range("K58").Value = 1146.4445
$ws.Range("M58").Value = -943.4445000000001
$ws.Range("H132").Value = 2737.818
$ws.Range("I132").Value = 2811.6333
$ws.Range("J132").Value = 1999.6666
$ws.Range("K132").Value = 8434.8999
$ws.Range("L132").Value = 5998.9998
$ws.Range("M132").Value = -5904.8999
$ws.Range("N132").Value = -11058.9998
$ws.Range("H136").Value = 1987.7858
$ws.Range("I136").Value = 1146.4445
$ws.Range("K136").Value = 3439.3335
$ws.Range("M136").Value = -889.3335000000002
$ws.Range("H141").Value = 117998.86
$ws.Range("J141").Value = 132832
$ws.Range("L141").Value = 132832
$ws.Range("N141").Value = -143192
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1302.5
$ws.Range("J11").Value = 2495
$ws.Range("L11").Value = 7485
$ws.Range("N11").Value = -7765
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H121").Value = 364349.9
$ws.Range("J121").Value = 666966.7
$ws.Range("L121").Value = 2000900.1
$ws.Range("N121").Value = -2003520.1
$ws.Range("H129").Value = 10102008
$ws.Range("I129").Value = 18519010
$ws.Range("K129").Value = 55557030
$ws.Range("M129").Value = -55552030
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6728.9
$ws.Range("I102").Value = 2223.625
$ws.Range("K102").Value = 2223.625
$ws.Range("M102").Value = -601.625
$ws.Range("H132").Value = 2247.7917
$ws.Range("I132").Value = 2287.4736
$ws.Range("K132").Value = 6862.4208
$ws.Range("M132").Value = -4332.4208
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2849.9092
$ws.Range("I22").Value = 2293.625
$ws.Range("J22").Value = 4333.3335
$ws.Range("K22").Value = 2293.625
$ws.Range("L22").Value = 4333.3335
$ws.Range("M22").Value = -1998.625
$ws.Range("N22").Value = -4923.3335
$ws.Range("H27").Value = 2849.9092
$ws.Range("I27").Value = 2293.625
$ws.Range("J27").Value = 4333.3335
$ws.Range("K27").Value = 2293.625
$ws.Range("L27").Value = 4333.3335
$ws.Range("M27").Value = -2186.625
$ws.Range("N27").Value = -4547.3335
$ws.Range("H55").Value = 927.9048
$ws.Range("I55").Value = 557
$ws.Range("J55").Value = 1335.9
$ws.Range("K55").Value = 557
$ws.Range("L55").Value = 1335.9
$ws.Range("M55").Value = -384
$ws.Range("N55").Value = -1681.9
$ws.Range("H61").Value = 2472
$ws.Range("I61").Value = 2306.125
$ws.Range("J61").Value = 3799
$ws.Range("K61").Value = 2306.125
$ws.Range("L61").Value = 3799
$ws.Range("M61").Value = -2104.125
$ws.Range("N61").Value = -4203
$ws.Range("H82").Value = 1956.3572
$ws.Range("J82").Value = 1677
$ws.Range("L82").Value = 1677
$ws.Range("N82").Value = -2399
$ws.Range("H85").Value = 1956.3572
$ws.Range("J85").Value = 1677
$ws.Range("L85").Value = 1677
$ws.Range("N85").Value = -4173
$ws.Range("H93").Value = 661141.7
$ws.Range("J93").Value = 2530051.8
$ws.Range("L93").Value = 2530051.8
$ws.Range("N93").Value = -2532547.8
$ws.Range("H113").Value = 2472
$ws.Range("I113").Value = 2306.125
$ws.Range("J113").Value = 3799
$ws.Range("K113").Value = 2306.125
$ws.Range("L113").Value = 3799
$ws.Range("M113").Value = -136.125
$ws.Range("N113").Value = -8139
$ws.Range("H132").Value = 7177.5713
$ws.Range("I132").Value = 2963.8333
$ws.Range("K132").Value = 8891.499899999999
$ws.Range("M132").Value = -6361.499899999999
$ws.Range("H135").Value = 75214.164
$ws.Range("J135").Value = 75214.164
$ws.Range("L135").Value = 75214.164
$ws.Range("N135").Value = -85354.164
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 50035000
$ws.Range("J135").Value = 50035000
$ws.Range("L135").Value = 50035000
$ws.Range("N135").Value = -50045140
